$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "('Beast of Burden', ['{6}', 'Artifact Creature — Golem', 'Beast of Burden’s power and toughness are each equal to the number of creatures on the battlefield.', '*/*'])"

$ws.Range("A3:A11").ClearContents()
